$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$wb.Activate()
$ws.Activate()

# Update the D-column prices (step 1 and 2 price list fix)
$ws.Range("D28").Value = 889.514
$ws.Range("D29").Value = 1016.592
$ws.Range("D30").Value = 1163.322
$ws.Range("D31").Value = 1360.888
$ws.Range("D32").Value = 1488.922
$ws.Range("D33").Value = 1774.242
$ws.Range("D34").Value = 1817.394
$ws.Range("D35").Value = 2027.423
$ws.Range("D36").Value = 2455.156
$ws.Range("D37").Value = 2771.646
$ws.Range("D38").Value = 3236.783
$ws.Range("D39").Value = 3572.451
$ws.Range("D40").Value = 3884.139
$ws.Range("D41").Value = 4363.662
$ws.Range("D42").Value = 4656.172
$ws.Range("D43").Value = 5078.153
$ws.Range("D49").Value = 2853.166
$ws.Range("D50").Value = 2867.551
$ws.Range("D51").Value = 3519.701
$ws.Range("D52").Value = 3764.262
$ws.Range("D53").Value = 4009.287
$ws.Range("D54").Value = 4512.31
$ws.Range("D55").Value = 5250.776
$ws.Range("D56").Value = 5878.954
$ws.Range("D57").Value = 6722.914
$ws.Range("D58").Value = 7576.465
$ws.Range("D59").Value = 8223.814
$ws.Range("D60").Value = 9254.794
$ws.Range("D61").Value = 10209.04
$ws.Range("D62").Value = 10813.235
$ws.Range("D68").Value = 1151.335
$ws.Range("D69").Value = 1184.423
$ws.Range("D70").Value = 1340.267
$ws.Range("D71").Value = 1561.329
$ws.Range("D72").Value = 1805.888
$ws.Range("D73").Value = 2093.598
$ws.Range("D74").Value = 2344.868
$ws.Range("D75").Value = 2483.928
$ws.Range("D76").Value = 2925.088
$ws.Range("D77").Value = 3390.232
$ws.Range("D78").Value = 3777.201
$ws.Range("D79").Value = 4233.225
$ws.Range("D80").Value = 4723.304
$ws.Range("D81").Value = 5002.86
$ws.Range("D82").Value = 5423.408
$ws.Range("D83").Value = 5922.106
$ws.Range("D89").Value = 492.383
$ws.Range("D90").Value = 591.254
$ws.Range("D91").Value = 639.205
$ws.Range("D92").Value = 793.128
$ws.Range("D93").Value = 865.538
$ws.Range("D94").Value = 932.194
$ws.Range("D95").Value = 1022.822
$ws.Range("D96").Value = 1305.265
$ws.Range("D97").Value = 1419.386
$ws.Range("D98").Value = 1663.943
$ws.Range("D99").Value = 1879.729
$ws.Range("D100").Value = 2023.588
$ws.Range("D106").Value = 4560.261
$ws.Range("D107").Value = 4847.978
$ws.Range("D108").Value = 5130.898
$ws.Range("D109").Value = 5610.424
$ws.Range("D110").Value = 5643.988
$ws.Range("D111").Value = 6650.984
$ws.Range("D112").Value = 7873.763
$ws.Range("D113").Value = 8727.319
$ws.Range("D114").Value = 9542.513000000001
$ws.Range("D115").Value = 10535.123
$ws.Range("D116").Value = 11796.261
$ws.Range("D117").Value = 12462.806
$ws.Range("D118").Value = 13906.164
$ws.Range("D119").Value = 14745.332
